# Primitive layout of the application
#
# Moves the hidden "_GoBack" bookmark paragraph from just before the
# "Add funds to account" section up to just after the Login section's
# trailing blank paragraph, and inserts a new blue "Register / Login
# screen" explanatory paragraph (plus an extra blank paragraph) in its
# place.

$d = $word.ActiveDocument

function Find-ParaIndex($substr, $occurrence) {
    $count = 0
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text.Contains($substr)) {
            $count++
            if ($count -eq $occurrence) {
                return $i
            }
        }
    }
    return -1
}

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Step 1: strip the _GoBack bookmark out of its old paragraph -----
# That paragraph sits right after the "...time remaining / parking has
# expired..." paragraph, just before "Add funds to account". Doing this
# edit first (it's further down the document) means it cannot disturb
# the paragraph indices we still need to resolve for step 2.
$notifIdx = Find-ParaIndex "they have x amount of time remaining" 1
$oldBookmarkParaIdx = $notifIdx + 1
$oldBookmarkPara = $d.Paragraphs.Item($oldBookmarkParaIdx)
$oldBookmarkPara.Range.InsertXML("<w:p $wNs/>")

# --- Step 2: insert the new paragraphs after the Login section -------
# Right after the second "...the user is 'logged in'" paragraph (end of
# the Login section, just before the "Logged In" heading) there is a
# single blank paragraph. Replace that blank paragraph with: the
# original blank paragraph + a new blank paragraph + the new blue
# explanatory paragraph + the relocated _GoBack bookmark paragraph.
$loginIdx = Find-ParaIndex "successful then the user is" 2
$blankParaIdx = $loginIdx + 1
$blankPara = $d.Paragraphs.Item($blankParaIdx)

$colorRpr = '<w:rPr><w:color w:val="548DD4" w:themeColor="text2" w:themeTint="99"/></w:rPr>'

$newXml = "<w:p $wNs/>" + `
    "<w:p $wNs/>" + `
    "<w:p $wNs><w:pPr>$colorRpr</w:pPr>" + `
    "<w:r>$colorRpr<w:t>The Register / Login Screen will be the same initial screen the user sees</w:t></w:r>" + `
    "<w:r>$colorRpr<w:t>. A fluid transition between the options would be ideal (</w:t></w:r>" + `
    "<w:proofErr w:type=`"spellStart`"/>" + `
    "<w:r>$colorRpr<w:t>tresorit</w:t></w:r>" + `
    "<w:proofErr w:type=`"spellEnd`"/>" + `
    "<w:r>$colorRpr<w:t xml:space=`"preserve`"> application). </w:t></w:r>" + `
    "</w:p>" + `
    "<w:p $wNs><w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/></w:p>"

$blankPara.Range.InsertXML($newXml)
